# Merge the <id>...</id> runs (currently split into three runs:
# "<id>", the bare id text, and "</id>") into a single run per
# occurrence, keeping the Courier New / color 7f6000 / sz 18 formatting
# that the opening "<id>" (and closing "</id>") run already carries.
#
# This mirrors what Word itself does when the three adjacent runs are
# retyped as one contiguous piece of text: the resulting single run
# inherits the character formatting of the first run in the replaced
# range, and the separate runs collapse away.

$d = $word.ActiveDocument

$ids = @("p091r_1", "p091r_2")

foreach ($id in $ids) {
    $searchText = "<id>" + $id + "</id>"

    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        # Setting .Text to the exact same string Word already displays is a
        # no-op, so round-trip through a temporary placeholder string to
        # force Word to actually collapse the matched range into a single
        # run (using the first run's formatting) before writing back the
        # real text.
        $rng.Text = "__TMP_MERGE_PLACEHOLDER__"
        $rng.Text = $searchText
    }
}
